# Update factsheets with text edits from COMM
#
# The source data pipeline started emitting the "No. of 990 Filers w/ Gov
# Grants" counts as text (to match the rest of the already-text formatted
# columns) instead of numbers, added a couple of "0.00%"/"$0" placeholder
# edits for two zero-filer counties, and appended a "Total" row to the
# County sheet that was previously missing it.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $Range,
        [string]$Text
    )
    # Force the cell to stay a literal text value even when the text looks
    # numeric (e.g. "503"), matching how the workbook now stores these
    # figures as text strings everywhere else on these sheets. Re-apply the
    # "Normal" style afterwards so we don't leave a stray text number format
    # behind on the cell.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Overall sheet
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall.Range("A2") "503"

# ---------------------------------------------------------------------
# County sheet
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

Set-TextValue $wsCounty.Range("B2")  "21"
Set-TextValue $wsCounty.Range("B3")  "8"
Set-TextValue $wsCounty.Range("B4")  "265"
Set-TextValue $wsCounty.Range("B5")  "14"
Set-TextValue $wsCounty.Range("B6")  "13"
Set-TextValue $wsCounty.Range("B7")  "1"
Set-TextValue $wsCounty.Range("B8")  "3"
Set-TextValue $wsCounty.Range("B9")  "2"
Set-TextValue $wsCounty.Range("B10") "5"
Set-TextValue $wsCounty.Range("B11") "2"
Set-TextValue $wsCounty.Range("B12") "5"
Set-TextValue $wsCounty.Range("B13") "3"
Set-TextValue $wsCounty.Range("B14") "6"
Set-TextValue $wsCounty.Range("B15") "148"
Set-TextValue $wsCounty.Range("B16") "7"

# Eureka County (row 17) and Lander County (row 18) switch from raw "0"
# placeholders to formatted percent/currency placeholders.
Set-TextValue $wsCounty.Range("B17") "0.00%"
Set-TextValue $wsCounty.Range("C17") "`$0"
Set-TextValue $wsCounty.Range("D17") "0.00%"
Set-TextValue $wsCounty.Range("E17") "0.00%"
Set-TextValue $wsCounty.Range("F17") "0.00%"

Set-TextValue $wsCounty.Range("B18") "0.00%"
Set-TextValue $wsCounty.Range("C18") "`$0"
Set-TextValue $wsCounty.Range("D18") "0.00%"
Set-TextValue $wsCounty.Range("E18") "0.00%"
Set-TextValue $wsCounty.Range("F18") "0.00%"

# New Total row (row 19), matching the Total rows on the other sheets.
Set-TextValue $wsCounty.Range("A19") "Total"
Set-TextValue $wsCounty.Range("B19") "503"
Set-TextValue $wsCounty.Range("C19") "`$1,252,051,703"
Set-TextValue $wsCounty.Range("D19") "7.25%"
Set-TextValue $wsCounty.Range("E19") "-23.82%"
Set-TextValue $wsCounty.Range("F19") "73.76%"

# ---------------------------------------------------------------------
# Congressional District sheet
# ---------------------------------------------------------------------
$wsCD = $wb.Worksheets.Item("Congressional District")

Set-TextValue $wsCD.Range("B2") "96"
Set-TextValue $wsCD.Range("B3") "228"
Set-TextValue $wsCD.Range("B4") "89"
Set-TextValue $wsCD.Range("B5") "90"
Set-TextValue $wsCD.Range("B6") "503"

# ---------------------------------------------------------------------
# Size sheet
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")

Set-TextValue $wsSize.Range("B2") "138"
Set-TextValue $wsSize.Range("B3") "169"
Set-TextValue $wsSize.Range("B4") "79"
Set-TextValue $wsSize.Range("B5") "39"
Set-TextValue $wsSize.Range("B6") "62"
Set-TextValue $wsSize.Range("B7") "16"
Set-TextValue $wsSize.Range("B8") "503"

# ---------------------------------------------------------------------
# Subsector sheet
# ---------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")

Set-TextValue $wsSubsector.Range("B2")  "41"
Set-TextValue $wsSubsector.Range("B3")  "72"
Set-TextValue $wsSubsector.Range("B4")  "35"
Set-TextValue $wsSubsector.Range("B5")  "72"
Set-TextValue $wsSubsector.Range("B6")  "7"
Set-TextValue $wsSubsector.Range("B7")  "147"
Set-TextValue $wsSubsector.Range("B8")  "3"
Set-TextValue $wsSubsector.Range("B9")  "2"
Set-TextValue $wsSubsector.Range("B10") "41"
Set-TextValue $wsSubsector.Range("B11") "4"
Set-TextValue $wsSubsector.Range("B12") "75"
Set-TextValue $wsSubsector.Range("B13") "4"
Set-TextValue $wsSubsector.Range("B14") "503"
